$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 785
$ws.Range("B3").Value = 49
$ws.Range("B4").Value = 177
$ws.Range("B5").Value = 446
$ws.Range("B6").Value = 123
$ws.Range("B7").Value = 237

$ws.Range("R14").Select()
